$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.887.50"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.88"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.63"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5045"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2574"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06396"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07777"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.658.36"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.286"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5442"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7871"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.00"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.952.35"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.22"
$ws.Range("E19").Value = "  -3.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.397"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.968"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.009"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.867"
$ws.Range("E24").Value = "  -4.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.07"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1143"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.877"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.244"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04979"
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.275"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.200"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.528"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.375"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8933"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.604"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.138.81"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5553"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01558"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.006"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.702"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8177"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.97"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +8.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.780.12"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4526"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.26"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05077"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09522"
$ws.Range("E51").Value = "  +2.52%  "
